# New submission synced: 2026-02-08 19:48:45
# Sheet "JSS 3B" gets a new response row (row 5) appended, and the
# previously mis-typed "Admission No" value in C4 is corrected to a
# genuine number.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3B")

# Fix C4: was stored as text "43", should be the number 43.
$ws.Range("C4").Value = 43

# Append the new submission as row 5.
$ws.Range("A5").Value = "2026-02-08 19:48:45"
$ws.Range("B5").Value = "Ibrahim Abubakar Shettima "
$ws.Range("C5").Value = "Number 3"
$ws.Range("D5").Value = 7
